$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8708921670913696
$ws.Range("B1").Value = 6.209417343139648
$ws.Range("C1").Value = 2.865254163742065
$ws.Range("D1").Value = 1.957189917564392
$ws.Range("E1").Value = 1.862564325332642
